$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add the external workbook reference (interventions-over-time vectors are
#        written out to global variables in COVID19_App_Data(9).xlsx) ---------------
$extCell = $ws.Range("D1")
$extCell.Formula = "='[COVID19_App_Data(9).xlsx]COVID19_App_Data(9)'!A1"
$extCell.EntireColumn.Delete()

# --- 2. New legend rows 20-29: *_xxx_vector column names (column A) ----------------
$ws.Range("A20").Value = "*_si_vector"
$ws.Range("A21").Value = "*_sd_vector"
$ws.Range("A22").Value = "*_scr_vector"
$ws.Range("A23").Value = "*_hw_vector"
$ws.Range("A24").Value = "*_wah_vector"
$ws.Range("A25").Value = "*_sc_vector"
$ws.Range("A26").Value = "*_tb_vector"
$ws.Range("A27").Value = "*_cte_vector"
$ws.Range("A28").Value = "*_q_vector"
$ws.Range("A29").Value = "*_vc_vector"

# --- 3. New legend rows 20-29: descriptions (column B) ------------------------------
$ws.Range("B20").Value = "Coverage over time (self-isolation)"
$ws.Range("B21").Value = "Coverage over time (social distancing)"
$ws.Range("B22").Value = "Coverage over time (screening)"
$ws.Range("B23").Value = "Coverage/efficacy over time (handwashing)"
$ws.Range("B24").Value = "Coverage over time (working at home)"
$ws.Range("B25").Value = "Coverage over time (school closure)"
$ws.Range("B26").Value = "Coverage over time (travel ban)"
$ws.Range("B27").Value = "Coverage over time (shielding the elderly)"
$ws.Range("B28").Value = "Coverage over time (household isolation)"
$ws.Range("B29").Value = "Coverage over time (vaccination)"

# --- 4. Footnote moves from A21 to A31 and its wording changes ---------------------
$ws.Range("A31").Value = "*: either for baseline or hypothetical scenario"

# --- 5. Scroll / selection bookkeeping, matching the author's saved view -----------
$ws.Range("A32").Select()
